$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 900
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 900
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 900
$ws.Range("N12").Value = -1240
$ws.Range("M12").ClearContents()
# Row 15
$ws.Range("H15").Value = 1236.55
$ws.Range("I15").Value = 1236.55
$ws.Range("K15").Value = 3709.65
$ws.Range("M15").Value = -3540.65
# Row 51
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
# Row 55
$ws.Range("H55").Value = 133
$ws.Range("I55").Value = 133
$ws.Range("K55").Value = 133
$ws.Range("M55").Value = 81
# Row 80
$ws.Range("H80").Value = 498.16666
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 498.16666
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 1494.49998
$ws.Range("N80").Value = -3490.49998
$ws.Range("M80").ClearContents()
# Row 83
$ws.Range("H83").Value = 498.16666
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 498.16666
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 4483.49994
$ws.Range("N83").Value = -14467.49994
$ws.Range("M83").ClearContents()
# Row 106
$ws.Range("H106").Value = 5850.4
$ws.Range("I106").Value = 5313
$ws.Range("K106").Value = 5313
$ws.Range("M106").Value = -4682
# Row 112
$ws.Range("H112").Value = 1649.5714
$ws.Range("J112").Value = 1796.7222
$ws.Range("L112").Value = 5390.1666
$ws.Range("N112").Value = -7606.1666
# Row 137
$ws.Range("H137").Value = 1708.9667
$ws.Range("I137").Value = 1255.5714
$ws.Range("J137").Value = 2766.889
$ws.Range("K137").Value = 3766.7142
$ws.Range("L137").Value = 8300.667000000001
$ws.Range("M137").Value = -1216.7142
$ws.Range("N137").Value = -13400.667
# Row 138
$ws.Range("H138").Value = 2652.878
$ws.Range("I138").Value = 1834.2727
$ws.Range("J138").Value = 2953.0334
$ws.Range("K138").Value = 5502.8181
$ws.Range("L138").Value = 8859.100199999999
$ws.Range("M138").Value = -362.8181000000004
$ws.Range("N138").Value = -19139.1002

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4901.02
$ws.Range("I32").Value = 3740
$ws.Range("J32").Value = 13415.167
$ws.Range("K32").Value = 3740
$ws.Range("L32").Value = 13415.167
$ws.Range("M32").Value = -3453
$ws.Range("N32").Value = -13989.167

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 694.625
$ws.Range("I22").Value = 651.1429000000001
$ws.Range("K22").Value = 651.1429000000001
$ws.Range("M22").Value = -478.1429000000001
# Row 134
$ws.Range("H134").Value = 1939.4445
$ws.Range("I134").Value = 1819.625
$ws.Range("K134").Value = 5458.875
$ws.Range("M134").Value = -2923.875

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 6
$ws.Range("H6").Value = 1665.3334
$ws.Range("I6").Value = 1000
$ws.Range("J6").Value = 1998
$ws.Range("K6").Value = 1000
$ws.Range("L6").Value = 1998
$ws.Range("M6").Value = -887
$ws.Range("N6").Value = -2224
# Row 7
$ws.Range("H7").Value = 480.5
$ws.Range("I7").Value = 357.33334
$ws.Range("K7").Value = 357.33334
$ws.Range("M7").Value = -244.33334
# Row 31
$ws.Range("H31").Value = 1000
$ws.Range("I31").Value = 1000
$ws.Range("K31").Value = 1000
$ws.Range("M31").Value = -705
# Row 34
$ws.Range("H34").Value = 1000
$ws.Range("I34").Value = 1000
$ws.Range("K34").Value = 1000
$ws.Range("M34").Value = -798
# Row 86
$ws.Range("H86").Value = 9743.416999999999
$ws.Range("I86").Value = 10873.375
$ws.Range("K86").Value = 10873.375
$ws.Range("M86").Value = -9750.375
# Row 89
$ws.Range("H89").Value = 9743.416999999999
$ws.Range("I89").Value = 10873.375
$ws.Range("K89").Value = 54366.875
$ws.Range("M89").Value = -48750.875
# Row 122
$ws.Range("H122").Value = 3031.7693
$ws.Range("I122").Value = 2255.2856
$ws.Range("K122").Value = 6765.8568
$ws.Range("M122").Value = -4315.8568

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 51
$ws.Range("H51").Value = 2000
$ws.Range("I51").Value = 2000
$ws.Range("K51").Value = 6000
$ws.Range("M51").Value = -5540
# Row 131
$ws.Range("H131").Value = 1499.3334
$ws.Range("J131").Value = 2498
$ws.Range("L131").Value = 7494
$ws.Range("N131").Value = -17574

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 52
$ws.Range("H52").Value = 20033
$ws.Range("J52").Value = 20033
$ws.Range("L52").Value = 20033
$ws.Range("N52").Value = -20551
# Row 70
$ws.Range("H70").Value = 5500
$ws.Range("I70").Value = 5000
$ws.Range("K70").Value = 5000
$ws.Range("M70").Value = -4730
# Row 73
$ws.Range("H73").Value = 5500
$ws.Range("I73").Value = 5000
$ws.Range("K73").Value = 5000
$ws.Range("M73").Value = -4064
# Row 80
$ws.Range("H80").Value = 5113.4614
$ws.Range("I80").Value = 4496.25
$ws.Range("J80").Value = 6101
$ws.Range("K80").Value = 4496.25
$ws.Range("L80").Value = 6101
$ws.Range("M80").Value = -3498.25
$ws.Range("N80").Value = -8097
# Row 83
$ws.Range("H83").Value = 5113.4614
$ws.Range("I83").Value = 4496.25
$ws.Range("J83").Value = 6101
$ws.Range("K83").Value = 22481.25
$ws.Range("L83").Value = 30505
$ws.Range("M83").Value = -17489.25
$ws.Range("N83").Value = -40489
# Row 132
$ws.Range("H132").Value = 2649.476
$ws.Range("I132").Value = 1895.7693
$ws.Range("K132").Value = 5687.3079
$ws.Range("M132").Value = -3157.3079

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 46
$ws.Range("H46").Value = 64955
$ws.Range("J46").Value = 64955
$ws.Range("L46").Value = 64955
$ws.Range("N46").Value = -65417
# Row 125
$ws.Range("H125").Value = 78775
$ws.Range("J125").Value = 78775
$ws.Range("L125").Value = 78775
$ws.Range("N125").Value = -88615
# Row 134
$ws.Range("H134").Value = 64955
$ws.Range("J134").Value = 64955
$ws.Range("L134").Value = 194865
$ws.Range("N134").Value = -199935
# Row 136
$ws.Range("H136").Value = 2051.2856
$ws.Range("I136").Value = 1691.3334
$ws.Range("K136").Value = 5074.0002
$ws.Range("M136").Value = -2524.0002
# Row 138
$ws.Range("H138").Value = 114450
$ws.Range("J138").Value = 114450
$ws.Range("L138").Value = 114450
$ws.Range("N138").Value = -124730
